$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26's date cell (A26) loses the "last row" date-only format and takes on
# the regular datetime format used by all the other data rows.
$ws.Range("A26").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 27: raw and clean SSA data for June 26th.
$ws.Range("A27").Value = 44008
$ws.Range("B27").Value = 208392
$ws.Range("C27").Value = 267288
$ws.Range("D27").Value = 66440
$ws.Range("E27").Value = 25779
$ws.Range("F27").Value = 31.37

# The new last row takes on the date-only format that row 26 used to have.
$ws.Range("A27").NumberFormat = "YYYY-MM-DD"
